# Apply odds updates to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 3.75
$ws.Range("H3").Value = 2.26
$ws.Range("Q3").Value = 1.82

# Row 5
$ws.Range("F5").Value = 2.9
$ws.Range("G5").Value = 3.65
$ws.Range("H5").Value = 2.26
$ws.Range("I5").Value = 2.74
$ws.Range("J5").Value = 3.45
$ws.Range("K5").Value = 4.1
$ws.Range("P5").Value = 2.06

# Row 8
$ws.Range("F8").Value = 2.28
$ws.Range("Q8").Value = 3

# Row 11
$ws.Range("F11").Value = 2.28
$ws.Range("G11").Value = 2.48
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 4.4

# Row 12
$ws.Range("F12").Value = 1.5
$ws.Range("H12").Value = 8.800000000000001
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 3.95
$ws.Range("K12").Value = 4.4
$ws.Range("P12").Value = 1.67
$ws.Range("Q12").Value = 2.28

# Row 13
$ws.Range("F13").Value = 1.41
$ws.Range("G13").Value = 1.51
$ws.Range("H13").Value = 8.6
$ws.Range("I13").Value = 12
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 5.3
$ws.Range("P13").Value = 2
$ws.Range("Q13").Value = 1.81
